$wb = $excel.ActiveWorkbook

$wsSystem = $wb.Worksheets.Item("system")
$wsSpecies = $wb.Worksheets.Item("species")

# --- system sheet data edits ---
# A2: 400 -> 250
$wsSystem.Range("A2").Value = 250
# E2 formula: 990000/1000000 -> 5000/1000000
$wsSystem.Range("E2").Formula = "=5000/1000000"

# --- species sheet data edits ---
# n-Hexane row (row 2): Inlet Mole Frac 0.3 -> 0.25, Relative Volatility 12 -> 8.9
$wsSpecies.Range("C2").Value = 0.25
$wsSpecies.Range("D2").Value = 8.9
# n-Heptane row (row 3): Relative Volatility 3.2 -> 5.7
$wsSpecies.Range("D3").Value = 5.7
# n-Octane row (row 4): Inlet Mole Frac 0.4 -> 0.45, Relative Volatility 1 -> 3.2
$wsSpecies.Range("C4").Value = 0.45
$wsSpecies.Range("D4").Value = 3.2

# --- selection / active sheet changes ---
# Selection on the "system" sheet moves to C2 and it is no longer the active tab
$wsSystem.Range("C2").Select()
# The "species" sheet becomes the active tab with selection C6
$wsSpecies.Activate()
$wsSpecies.Range("C6").Select()

# --- workbook window geometry ---
$excel.Width = 29040
$excel.Height = 15720
$excel.Left = 28680
$excel.Top = -2190
